# Insert a new weekly record into the "Poroto granado" sheet.
# The new record is inserted right after the existing row 119 (becoming the
# new row 120), which pushes the old rows 120-126 down by one (to 121-127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 120, shifting existing rows 120+ down.
$ws.Rows.Item(120).Insert()

# Fill in the new row 120 with the new weekly data, copying the constant /
# repeated columns from the row above (row 119) and setting the unique
# values per the diff.
$ws.Range("A120").Value = 8
$ws.Range("B120").Value = "Terminal La Palmera de La Serena"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 45008
$ws.Range("D120").NumberFormat = $ws.Range("D119").NumberFormat
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = 100112030
$ws.Range("G120").Value = "Poroto granado"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 400
$ws.Range("K120").Value = 34000
$ws.Range("L120").Value = 35000
$ws.Range("M120").Value = 34500
$ws.Range("N120").Value = "$/malla 25 kilos"
$ws.Range("O120").Value = "Provincia del Elquí"
$ws.Range("P120").Value = 1380
$ws.Range("Q120").Value = 25
$ws.Range("R120").Value = "Hortaliza"
